# Update the "取得日時" (acquired datetime) column for all data rows
# on the active sheet ("ランサーズ") from 2025-09-26 06:29:13 to
# 2025-09-26 06:35:00, matching a fresh scrape append timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "2025-09-26 06:29:13"
$newValue = "2025-09-26 06:35:00"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
